# Cryptology7-Elliptic-Curve.pptx -- "changes to course end"
#
# Two content edits are applied, found by locating their shapes/text and
# editing via TextRange.Characters(start,length) so existing run
# formatting (size/color/hyperlinks/etc.) on untouched text is preserved
# and only the minimal affected span is rewritten.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 11: "key = bA mod p = B(aP) mod p" -> "... = b(aP) mod p"
# (capitalized "B(" becomes lowercase "b" + "(")
# ---------------------------------------------------------------------
$slide11 = $p.Slides.Item(11)
for ($i = 1; $i -le $slide11.Shapes.Count; $i++) {
    $shape = $slide11.Shapes.Item($i)
    if (-not $shape.HasTextFrame) { continue }
    $tr = $shape.TextFrame.TextRange
    $full = $tr.Text
    if ($full.IndexOf("B(aP) mod p") -ge 0) {
        $idx0 = $full.IndexOf("B(")
        $pos = $idx0 + 1
        # Rewrite just the "B" -> "b"; the trailing "(" stays where it is
        # but becomes its own run since the preceding run was edited.
        $tr.Characters($pos, 1).Text = "b"
        break
    }
}

# ---------------------------------------------------------------------
# Slide 13: "Expected to become predominant method soon"
#        -> "Is now predominant method"
# (kept as two runs: "Is now " + "predominant method")
# ---------------------------------------------------------------------
$slide13 = $p.Slides.Item(13)
for ($i = 1; $i -le $slide13.Shapes.Count; $i++) {
    $shape = $slide13.Shapes.Item($i)
    if (-not $shape.HasTextFrame) { continue }
    $tr = $shape.TextFrame.TextRange
    $full = $tr.Text
    $oldPhrase = "Expected to become predominant method soon"
    $idx0 = $full.IndexOf($oldPhrase)
    if ($idx0 -ge 0) {
        $pos = $idx0 + 1
        $prefixOld = "Expected to become "
        $prefixNew = "Is now "
        $midText = "predominant method"
        $suffixOld = " soon"

        # Replace the prefix in place; the untouched middle text then
        # naturally becomes its own run (no characters inside it change).
        $tr.Characters($pos, $prefixOld.Length).Text = $prefixNew

        $midPos = $pos + $prefixNew.Length
        $suffixPos = $midPos + $midText.Length

        # Drop the trailing " soon".
        $tr.Characters($suffixPos, $suffixOld.Length).Text = ""
        break
    }
}
